# Generate Report for Handback
# This script re-applies the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" info for the 7f297600-...md item
# (row 7) on both the zh-cn and de-de localization-status worksheets, because
# a new (non-latest) handback file was detected for it.

$wb = $excel.ActiveWorkbook

# Hyperlink-like font used throughout the workbook for the "file name" /
# "target file" columns (matches the existing custom "HyperLink" cell style).
$hyperlinkColor = 15570276  # RGB(0x64, 0x95, 0xED) == FF6495ED

function Update-LanguageSheet {
    param(
        [string]$SheetName,
        [string]$XlfFileName,
        [string]$HandbackDateTime,
        [string]$CurrentCommit,
        [string]$LatestCommit
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Column P ("Error Detail") needs to be wide enough to show the message.
    $ws.Columns.Item(16).ColumnWidth = 39.1666666667

    $mdFileName = "7f297600-d8a3-498a-bc98-366c371e6dd7.md"

    # I7 = "Latest Target File": hyperlink to the actual (current) handback
    # commit, same display text as the "Source File Name" handback link.
    $currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$CurrentCommit/e2e/$mdFileName"
    $ws.Hyperlinks.Add($ws.Range("I7"), $currentUrl, [Type]::Missing, [Type]::Missing, $mdFileName) | Out-Null
    $ws.Range("I7").Font.Color = $hyperlinkColor
    $ws.Range("I7").Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle

    # J7 = "Latest Handback File": the xlf file that was handed back.
    $ws.Range("J7").Value = $XlfFileName

    # K7 = "Latest Handback DateTime": when that handback happened.
    $ws.Range("K7").Value = $HandbackDateTime

    # P7 = "Error Detail": the handback file isn't the latest one available.
    $latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$LatestCommit/e2e/$mdFileName"
    $errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
    $ws.Range("P7").Value = $errorDetail
}

Update-LanguageSheet "zh-cn" `
    "7f297600-d8a3-498a-bc98-366c371e6dd7.75a877359c12ae9e927da03389b2549b4aab4bb1.zh-cn.xlf" `
    "2016-08-31 00:46:31" `
    "a8b942eba3f8b111a0ffecdf557321f2d610c5c7" `
    "b0fa40f6f750d2706d121d55f02ec9aed5f4ef82"

Update-LanguageSheet "de-de" `
    "7f297600-d8a3-498a-bc98-366c371e6dd7.75a877359c12ae9e927da03389b2549b4aab4bb1.de-de.xlf" `
    "2016-08-31 00:46:38" `
    "a8b942eba3f8b111a0ffecdf557321f2d610c5c7" `
    "b0fa40f6f750d2706d121d55f02ec9aed5f4ef82"

Write-Output "Done updating handback report."
